$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data through row 182 (A1:V182). We need to append
# three new match rows (183-185), matching the formatting of the existing
# data rows: column A uses style index 1 (bold/bordered/centered "Indice"
# style) and column E uses style index 2 (datetime number format); the rest
# of the columns use the default style.

# Seed the formatting of the new rows by copying the last existing data row
# (182) and pasting only the formats into rows 183-185.
$ws.Range("A182:V182").Copy()
$ws.Range("A183:V183").PasteSpecial(-4122)
$ws.Range("A182:V182").Copy()
$ws.Range("A184:V184").PasteSpecial(-4122)
$ws.Range("A182:V182").Copy()
$ws.Range("A185:V185").PasteSpecial(-4122)

# --- Row 183: Bucaramanga 2 x 3 La Equidad ---
$ws.Range("A183").Value = 182
$ws.Range("B183").Value = "colombia"
$ws.Range("C183").Value = "primera-a"
$ws.Range("D183").Value = "'2023"
$ws.Range("E183").Value = 45224.95833333334
$ws.Range("F183").Value = "Bucaramanga"
$ws.Range("G183").Value = 2
$ws.Range("H183").Value = "La Equidad"
$ws.Range("I183").Value = 3
$ws.Range("J183").Value = 2.63
$ws.Range("K183").Value = "22/10/2023 22:42"
$ws.Range("L183").Value = 2.85
$ws.Range("M183").Value = "25/10/2023 22:55"
$ws.Range("N183").Value = 2.92
$ws.Range("O183").Value = "22/10/2023 22:42"
$ws.Range("P183").Value = 2.97
$ws.Range("Q183").Value = "25/10/2023 22:55"
$ws.Range("R183").Value = 3.12
$ws.Range("S183").Value = "22/10/2023 22:42"
$ws.Range("T183").Value = 2.89
$ws.Range("U183").Value = "25/10/2023 22:55"
$ws.Range("V183").Value = "https://www.betexplorer.com/football/colombia/primera-a/bucaramanga-la-equidad/WYDw1ePI/"

# --- Row 184: Pereira 0 x 2 Junior ---
$ws.Range("A184").Value = 183
$ws.Range("B184").Value = "colombia"
$ws.Range("C184").Value = "primera-a"
$ws.Range("D184").Value = "'2023"
$ws.Range("E184").Value = 45225.04861111111
$ws.Range("F184").Value = "Pereira"
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = "Junior"
$ws.Range("I184").Value = 2
$ws.Range("J184").Value = 2.35
$ws.Range("K184").Value = "23/10/2023 00:12"
$ws.Range("L184").Value = 2.78
$ws.Range("M184").Value = "26/10/2023 01:04"
$ws.Range("N184").Value = 2.97
$ws.Range("O184").Value = "23/10/2023 00:12"
$ws.Range("P184").Value = 3
$ws.Range("Q184").Value = "26/10/2023 01:01"
$ws.Range("R184").Value = 3.55
$ws.Range("S184").Value = "23/10/2023 00:12"
$ws.Range("T184").Value = 2.95
$ws.Range("U184").Value = "26/10/2023 01:01"
$ws.Range("V184").Value = "https://www.betexplorer.com/football/colombia/primera-a/dep-pereira-junior/O8AV2Hf6/"

# --- Row 185: Millonarios 1 x 1 Chico ---
$ws.Range("A185").Value = 184
$ws.Range("B185").Value = "colombia"
$ws.Range("C185").Value = "primera-a"
$ws.Range("D185").Value = "'2023"
$ws.Range("E185").Value = 45225.13888888889
$ws.Range("F185").Value = "Millonarios"
$ws.Range("G185").Value = 1
$ws.Range("H185").Value = "Chico"
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = 1.49
$ws.Range("K185").Value = "22/10/2023 22:42"
$ws.Range("L185").Value = 1.42
$ws.Range("M185").Value = "26/10/2023 03:11"
$ws.Range("N185").Value = 4.06
$ws.Range("O185").Value = "22/10/2023 22:42"
$ws.Range("P185").Value = 4.38
$ws.Range("Q185").Value = "26/10/2023 03:11"
$ws.Range("R185").Value = 7.75
$ws.Range("S185").Value = "22/10/2023 22:42"
$ws.Range("T185").Value = 9.44
$ws.Range("U185").Value = "26/10/2023 03:11"
$ws.Range("V185").Value = "https://www.betexplorer.com/football/colombia/primera-a/millonarios-chico/APCs0FvP/"

# Setting column D via a leading apostrophe (to force text, since "2023"
# looks numeric) marks the cell with a "quote prefix" style variant. Re-paste
# the plain formats for column D from the template row to drop that quote
# prefix marker, leaving the default (unstyled) cell format while keeping
# the stored value as text.
$ws.Range("D182").Copy()
$ws.Range("D183").PasteSpecial(-4122)
$ws.Range("D182").Copy()
$ws.Range("D184").PasteSpecial(-4122)
$ws.Range("D182").Copy()
$ws.Range("D185").PasteSpecial(-4122)

$excel.CutCopyMode = $false
